$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Fix the header row: the column names were missing the underscore
# between "Output" and the trailing digit ---
$ws.Range("A1").Value = "Output_7"
$ws.Range("B1").Value = "Output_1"
$ws.Range("C1").Value = "Output_2"

# --- Add a new column D, reusing the same header text as column C, and
# fill every data row with the literal text value "6" (not the number 6).
# Building it through a formula and then freezing the result with
# Copy + PasteSpecial(values) keeps the cell a plain shared-string value
# instead of turning it into a number or leaving a live formula behind. ---
$ws.Range("D1").Value = "Output_2"

$xlPasteValues = -4163
$dataRows = 2..7
foreach ($r in $dataRows) {
    $cell = $ws.Cells.Item($r, 4)
    $cell.Formula = "=""6"""
    $cell.Copy()
    $cell.PasteSpecial($xlPasteValues)
}

$excel.CutCopyMode = $false
